$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.47694268276958951,
    -0.087913362913113247,
    -0.0089999993895677477,
    -0.011999999845025755,
    -0.0059999993962627229,
    -0.0059999993781545413,
    -0.019999999275297498,
    -0.019999999273843549,
    -0.0059999993759793924,
    -0.0059999993770532001,
    0.06603235845526001,
    -0.0059999993754744629,
    -0.005999999365508657,
    -0.011999999318904386,
    -0.0059999993605117652,
    -0.0059999993585497791,
    -0.0059999993558799147,
    -0.0089999993334215489,
    -0.093858351634332315,
    -0.0089999994002454287,
    -0.0089999993994398508,
    -0.0089999993989149374,
    -0.0089999993840681469,
    -0.041999999134939792,
    -0.041999999130128529,
    -0.0059999993761934434,
    -0.0059999993729484835,
    -0.0059999993588428779,
    -0.011999999305819742,
    -0.01999999924290341,
    -0.01499999927272988,
    -0.02099999922795881,
    -0.0059999993369705962
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
